# Add a "Save" column (H) to the s_vals sheet, matching the header style
# used by the existing columns (B1:G1) and filling in data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the neighboring header cell (G1) onto the new header
# cell (H1) so it reuses the same bold/border/centered style instead of
# creating a brand new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the header text and the new column's data values.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
